$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting existing rows 87-161 down to 88-162.
$ws.Rows("87:87").Insert()

# Populate the newly inserted row 87 with the new weekly record.
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = (Get-Date -Year 2021 -Month 10 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112044
$ws.Range("G87").Value = "Perejil"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 180
$ws.Range("K87").Value = 4500
$ws.Range("L87").Value = 4500
$ws.Range("M87").Value = 4500
$ws.Range("N87").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O87").Value = "Región Metropolitana"
$ws.Range("P87").Value = 1500
$ws.Range("Q87").Value = 3
$ws.Range("R87").Value = "Hortaliza"
